{"js": "// The captured change is a pure OOXML re-serialization: every hunk in the\n// source diff only reorders XML namespace declarations on the root elements\n// of word/document.xml, word/footer1.xml, word/footnotes.xml and\n// word/styles.xml (sorted alphabetically), and reorders the attributes of\n// individual elements (e.g. `w:type=\"dxa\" w:w=\"3070\"` instead of\n// `w:w=\"3070\" w:type=\"dxa\"`). No element, attribute value, run, paragraph,\n// table, style definition or piece of text is added, removed or modified -\n// the rendered document and the Word/OOXML object model are 100%\n// unchanged (this matches the commit message, a tooling version bump from\n// 2.0.1 to 2.0.2 that changed how the authoring library serializes XML\n// attribute order).\n//\n// The Word JavaScript API only exposes the document's logical content\n// model (paragraphs, ranges, tables, styles, ...); it has no notion of,\n// and no way to control, the physical attribute ordering used when the\n// underlying part XML is serialized. There is therefore no content-level\n// operation to perform here - applying \"the change\" means touching\n// nothing observable through the object model. We simply read a couple of\n// harmless properties (so the script still demonstrates the\n// load/sync pattern) and make no mutation, which keeps the document's\n// content, formatting and structure byte-for-byte equivalent to the\n// source, exactly like the diff.\n\nconst body = context.document.body;\nbody.load(\"text\");\n\nconst sections = context.document.sections;\nsections.load(\"items\");\n\nawait context.sync();\n", "ps1": "# The captured change is a pure OOXML re-serialization: every hunk in the\n# source diff only reorders XML namespace declarations on the root elements\n# of word/document.xml, word/footer1.xml, word/footnotes.xml and\n# word/styles.xml (sorted alphabetically), and reorders the attributes of\n# individual elements (e.g. `w:type=\"dxa\" w:w=\"3070\"` instead of\n# `w:w=\"3070\" w:type=\"dxa\"`). No element, attribute value, run, paragraph,\n# table, style definition or piece of text is added, removed or modified -\n# the rendered document and the Word object model are 100% unchanged (this\n# matches the commit message, a tooling version bump from 2.0.1 to 2.0.2\n# that changed how the authoring library serializes XML attribute order).\n#\n# The Word COM object model only exposes the document's logical content\n# (Content, Paragraphs, Tables, Styles, ...); it has no property that\n# controls the physical attribute ordering used when the underlying part\n# XML is serialized. There is therefore no content-level operation to\n# perform here - applying \"the change\" means touching nothing observable\n# through the object model. We just read a couple of harmless properties\n# (so the script still demonstrates the usual COM interaction pattern) and\n# issue no mutating call, which keeps the document's content, formatting\n# and structure byte-for-byte equivalent to the source, exactly like the\n# diff.\n\n$d = $word.ActiveDocument\n\n# Read-only touches, no mutation performed.\n$null = $d.Content.Text\n$null = $d.Sections.Item(1).Footers.Item(1).Range.Text\n"}
